# Africa weekly variant summary - update voi/voc table structure and refresh
# data up to 2021-10-14:
#  - inserted a new column F ("new sequences sampled and submitted in the
#    last 30 days"), pushing the old F/G/H (dates, days-since) one column
#    right to G/H/I
#  - reordered / replaced the data rows with the refreshed counts
#  - widened the used range from A1:H10 to A1:I10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Wipe the previous data area (A2:H10) so stale cells (e.g. a B column
#    that used to be empty for some rows) don't linger once the columns
#    shift around.
# ---------------------------------------------------------------------------
$ws.Range("A2:I10").ClearContents()

# ---------------------------------------------------------------------------
# 2. Header row (row 1) - E1/F1 wording changes, F1 is a brand-new column,
#    and the rest shift from F/G/H -> G/H/I.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "Variant (VOC,VOI,VUM)"
$ws.Cells.Item(1,2).Value = "Other names by which this variant may be known"
$ws.Cells.Item(1,3).Value = "Lineage/sub-lineages"
$ws.Cells.Item(1,4).Value = "Number of sequences"
$ws.Cells.Item(1,5).Value = "Sequences submitted in the last 30 days (data up to 2021-10-14)"
$ws.Cells.Item(1,6).Value = "new sequences sampled and submitted in the last 30 days (data up to 2021-10-14)"
$ws.Cells.Item(1,7).Value = "date of first sequence"
$ws.Cells.Item(1,8).Value = "date of last sequence"
$ws.Cells.Item(1,9).Value = "No of days since last sampling"

# ---------------------------------------------------------------------------
# 3. Data rows 2-10. Columns G/H hold dates that must stay plain text (not
#    be auto-parsed into date serials), so format those two columns as text
#    up front, then fill in row-by-row. Missing entries (some rows have no
#    "other names" / no "submitted in last 30 days" / no "new sequences")
#    are simply left blank - ClearContents() above already removed them.
# ---------------------------------------------------------------------------
$ws.Range("G2:H10").NumberFormat = "@"

$data = @(
  # row, A,                B,                 C,                                                  D,     E,    F,   G,            H
  @(2,  "Alpha",          "VOC-202012/01",    "B.1.1.7",                                           2564,  399,  1,  "2020-08-02", "2021-09-30", 16),
  @(3,  "Beta",           "VOC-202012/02",    "B.1.351, B.1.351.1, B.1.351.2",                      10219, 437,  1,  "2020-05-27", "2021-09-25", 21),
  @(4,  "Delta",          "VOC-21APR-02",     "B.1.617.2, AY.1, AY.3, AY.4, AY.5, AY.6, AY.7.1, AY.10, AY.11, AY.12, AY.13, AY.14, AY.15, AY.16, AY.17, AY.18, AY.19, AY.20, AY.21, AY.23, AY.24, AY.25", 12722, 3444, 52, "2020-09-10", "2021-09-20", 26),
  @(5,  "B.1.1.318",      "VUM-2021-06-04",   "B.1.1.318, AZ.1, AZ.2, AZ.5",                        672,   110,  $null, "2021-01-06", "2021-08-11", 66),
  @(6,  "C.1",            $null,              "C.1",                                               384,   $null,$null, "2020-04-12", "2021-06-28", 110),
  @(7,  "C.1.2",          "VUM-2021-09-01",   "C.1.2",                                              187,   48,   4,  "2021-05-11", "2021-09-17", 29),
  @(8,  "C.36.3",         "VUM-2021-06-16",   "C.36.3",                                             104,   2,    $null, "2020-04-26", "2021-06-16", 122),
  @(9,  "Eta",            "VUM-202102/03",    "B.1.525",                                            991,   197,  3,  "2020-03-28", "2021-09-17", 29),
  @(10, "A.23.1",         $null,              "A.23.1",                                             388,   26,   $null, "2020-10-21", "2021-08-22", 55)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r,1).Value = $row[1]
  if ($null -ne $row[2]) { $ws.Cells.Item($r,2).Value = $row[2] }
  $ws.Cells.Item($r,3).Value = $row[3]
  $ws.Cells.Item($r,4).Value = $row[4]
  if ($null -ne $row[5]) { $ws.Cells.Item($r,5).Value = $row[5] }
  if ($null -ne $row[6]) { $ws.Cells.Item($r,6).Value = $row[6] }
  $ws.Cells.Item($r,7).Value = $row[7]
  $ws.Cells.Item($r,8).Value = $row[8]
  $ws.Cells.Item($r,9).Value = $row[9]
}
